# Update workbook with new hydrological readings (atualizacao 16 nov 2020)
$wb = $excel.ActiveWorkbook

# ---- Sheet "Mensal" (first sheet): add monthly summary row for Nov/2020 ----
$wsMensal = $wb.Worksheets.Item(1)
$lastMensalRow = 13
$newMensalRow = $lastMensalRow + 1

# Copy formatting (incl. date number format) from the row above onto the new row
$wsMensal.Cells.Item($lastMensalRow, 1).Copy($wsMensal.Cells.Item($newMensalRow, 1))

$wsMensal.Cells.Item($newMensalRow, 1).Value = 44150
$wsMensal.Cells.Item($newMensalRow, 2).Value = 69.02
$wsMensal.Cells.Item($newMensalRow, 3).Value = 140.62
$wsMensal.Cells.Item($newMensalRow, 4).Value = -50.92

# ---- Sheet "Diario" (second sheet): append daily readings through 2020-11-15 ----
$wsDiario = $wb.Worksheets.Item(2)
$lastDiarioRow = 367

$dailyRows = @(
    @(44136, 70.31, 140.62, -50),
    @(44137, 69.23999999999999, 140.62, -50.76),
    @(44138, 69.66, 140.62, -50.47),
    @(44139, 69.95, 140.62, -50.26),
    @(44140, 67.93000000000001, 140.62, -51.7),
    @(44141, 67.36, 140.62, -52.1),
    @(44142, 66.19, 140.62, -52.93),
    @(44143, 65.41, 140.62, -53.49),
    @(44144, 68.23999999999999, 140.62, -51.47),
    @(44145, 66.95, 140.62, -52.39),
    @(44146, 67.65000000000001, 140.62, -51.9),
    @(44147, 70.76000000000001, 140.62, -49.68),
    @(44148, 70.52, 140.62, -49.85),
    @(44149, 71.98, 140.62, -48.81),
    @(44150, 73.11, 140.62, -48.01)
)

$row = $lastDiarioRow
foreach ($vals in $dailyRows) {
    $row = $row + 1
    # Copy formatting (incl. date number format) from the previous row onto the new row
    $wsDiario.Cells.Item($row - 1, 1).Copy($wsDiario.Cells.Item($row, 1))

    $wsDiario.Cells.Item($row, 1).Value = $vals[0]
    $wsDiario.Cells.Item($row, 2).Value = $vals[1]
    $wsDiario.Cells.Item($row, 3).Value = $vals[2]
    $wsDiario.Cells.Item($row, 4).Value = $vals[3]
}
